$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 102.5
$ws.Range("I5").Value = 102.5
$ws.Range("K5").Value = 102.5
$ws.Range("M5").Value = 12.5

$ws.Range("H40").Value = 4918.72
$ws.Range("I40").Value = 3373.9375
$ws.Range("J40").Value = 7665
$ws.Range("K40").Value = 3373.9375
$ws.Range("L40").Value = 7665
$ws.Range("M40").Value = -3198.9375
$ws.Range("N40").Value = -8015

$ws.Range("H93").Value = 17999
$ws.Range("J93").Value = 17999
$ws.Range("L93").Value = 17999
$ws.Range("N93").Value = -22991

$ws.Range("H123").Value = 72495.75
$ws.Range("J123").Value = 72495.75
$ws.Range("L123").Value = 72495.75
$ws.Range("N123").Value = -82295.75

$ws.Range("H132").Value = 2303.1785
$ws.Range("I132").Value = 1370.3914
$ws.Range("K132").Value = 4111.174199999999
$ws.Range("M132").Value = -1581.174199999999

$ws.Range("H138").Value = 4544.0293
$ws.Range("I138").Value = 2910.35
$ws.Range("J138").Value = 5224.729
$ws.Range("K138").Value = 8731.049999999999
$ws.Range("L138").Value = 15674.187
$ws.Range("M138").Value = -3591.049999999999
$ws.Range("N138").Value = -25954.187

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7050.5
$ws.Range("I45").Value = 4586.857
$ws.Range("J45").Value = 10499.6
$ws.Range("K45").Value = 4586.857
$ws.Range("L45").Value = 10499.6
$ws.Range("M45").Value = -4209.857
$ws.Range("N45").Value = -11253.6

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H92").Value = 16696667
$ws.Range("J92").Value = 45000
$ws.Range("L92").Value = 45000
$ws.Range("N92").Value = -49992

$ws.Range("H110").Value = 715548.1
$ws.Range("I110").Value = 715548.1
$ws.Range("K110").Value = 715548.1
$ws.Range("M110").Value = -713503.1

$ws.Range("H132").Value = 12635.429
$ws.Range("I132").Value = 6724
$ws.Range("K132").Value = 20172
$ws.Range("M132").Value = -17642

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 98439.8
$ws.Range("J35").Value = 98439.8
$ws.Range("L35").Value = 98439.8
$ws.Range("N35").Value = -99059.8

$ws.Range("H88").Value = 10816.5
$ws.Range("J88").Value = 10816.5
$ws.Range("L88").Value = 10816.5
$ws.Range("N88").Value = -11628.5

$ws.Range("H91").Value = 10816.5
$ws.Range("J91").Value = 10816.5
$ws.Range("L91").Value = 10816.5
$ws.Range("N91").Value = -13624.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 337035.84
$ws.Range("I58").Value = 1250699.9
$ws.Range("K58").Value = 1250699.9
$ws.Range("M58").Value = -1250496.9

$ws.Range("H95").Value = 13665.6
$ws.Range("J95").Value = 13665.6
$ws.Range("L95").Value = 13665.6
$ws.Range("N95").Value = -19157.6

$ws.Range("H99").Value = 5361.769
$ws.Range("I99").Value = 3968.8
$ws.Range("K99").Value = 3968.8
$ws.Range("M99").Value = -2470.8

$ws.Range("H126").Value = 5361.769
$ws.Range("I126").Value = 3968.8
$ws.Range("K126").Value = 11906.4
$ws.Range("M126").Value = -9436.400000000001

$ws.Range("H132").Value = 5966.1377
$ws.Range("I132").Value = 5526.8423
$ws.Range("K132").Value = 16580.5269
$ws.Range("M132").Value = -14050.5269

$ws.Range("H136").Value = 337035.84
$ws.Range("I136").Value = 1250699.9
$ws.Range("K136").Value = 3752099.7
$ws.Range("M136").Value = -3749549.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 17840626
$ws.Range("I4").Value = 1083754.8
$ws.Range("K4").Value = 3251264.4
$ws.Range("M4").Value = -3251152.4

$ws.Range("H5").Value = 1309
$ws.Range("I5").Value = 979.5
$ws.Range("K5").Value = 2938.5
$ws.Range("M5").Value = -2826.5

$ws.Range("H102").Value = 12256.5
$ws.Range("J102").Value = 15000
$ws.Range("L102").Value = 45000
$ws.Range("N102").Value = -49868

$ws.Range("H116").Value = 2606.2
$ws.Range("I116").Value = 1999.5
$ws.Range("J116").Value = 3010.6667
$ws.Range("K116").Value = 5998.5
$ws.Range("L116").Value = 9032.000100000001
$ws.Range("M116").Value = -2556.5
$ws.Range("N116").Value = -15916.0001

$ws.Range("H132").Value = 3640.4333
$ws.Range("J132").Value = 3939.5
$ws.Range("L132").Value = 35455.5
$ws.Range("N132").Value = -40515.5

$ws.Range("H135").Value = 1309
$ws.Range("I135").Value = 979.5
$ws.Range("K135").Value = 8815.5
$ws.Range("M135").Value = -6280.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5500000
$ws.Range("I11").Value = 1000000
$ws.Range("K11").Value = 1000000
$ws.Range("M11").Value = -999861

$ws.Range("H12").Value = 2002.6666
$ws.Range("J12").Value = 2002.6666
$ws.Range("L12").Value = 2002.6666
$ws.Range("M12").Value = -2282.6666

$ws.Range("H70").Value = 35719824
$ws.Range("I70").Value = 4920.0835
$ws.Range("J70").Value = 62506000
$ws.Range("K70").Value = 4920.0835
$ws.Range("L70").Value = 62506000
$ws.Range("M70").Value = -4650.0835
$ws.Range("N70").Value = -62506540

$ws.Range("H73").Value = 35719824
$ws.Range("I73").Value = 4920.0835
$ws.Range("J73").Value = 62506000
$ws.Range("K73").Value = 4920.0835
$ws.Range("L73").Value = 62506000
$ws.Range("M73").Value = -3984.0835
$ws.Range("N73").Value = -62507872

$ws.Range("H102").Value = 9424.5
$ws.Range("I102").Value = 9424.5
$ws.Range("K102").Value = 9424.5
$ws.Range("M102").Value = -7802.5

$ws.Range("H117").Value = 40000
$ws.Range("J117").Value = 40000
$ws.Range("L117").Value = 40000
$ws.Range("N117").Value = -46884

$ws.Range("H122").Value = 5956.846
$ws.Range("I122").Value = 5970.391
$ws.Range("K122").Value = 17911.173
$ws.Range("M122").Value = -15461.173

$ws.Range("H126").Value = 90912560
$ws.Range("J126").Value = 3702.75
$ws.Range("L126").Value = 11108.25
$ws.Range("N126").Value = -16048.25

$ws.Range("H132").Value = 1252552.1
$ws.Range("J132").Value = 2999.5
$ws.Range("L132").Value = 8998.5
$ws.Range("N132").Value = -14058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 15000

$ws.Range("H107").Value = 3998.5
$ws.Range("I107").Value = 3998.5
$ws.Range("K107").Value = 3998.5
$ws.Range("M107").Value = -2078.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 125000
$ws.Range("J92").Value = 125000
$ws.Range("L92").Value = 125000
$ws.Range("N92").Value = -129992
